$d = $word.ActiveDocument

# 1) Title: "GRAFISCHE RICHTLINIEN" -> "GRAFISCHE IDENTITÄT"
$d.Content.Find.Execute(" GRAFISCHE RICHTLINIEN", $true, $false, $false, $false, $false,
                         $true, 1, $false, " GRAFISCHE IDENTITÄT", 2)

# 2) "Graphic Identity Guidelines" -> "Grafische Identität Richtlinie"
$d.Content.Find.Execute("Graphic Identity Guidelines", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Grafische Identität Richtlinie", 2)

# 3) Long description paragraph (English -> German)
$d.Content.Find.Execute("These are recommended usage guidelines for maintaining a consistent design aesthetic for the SmartCash brand. A strong and consistent visual identity of our logo will help keep a consistent look, recognition and familiarity now and in the future. Standardization of colours will go a long way to enforce a reliable and positive impression to our identity in the blockchain space.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Diese Empfehlungen werden für die Aufrechterhaltung einer konsistenten Designästhetik für die Marke SmartCash empfohlen. Eine starke und konsequente visuelle Identität unseres Logos wird dazu beitragen, dass wir jetzt und in Zukunft ein konsistentes Aussehen, Anerkennung und Vertrautheit erhalten. Standardisierung von Farben wird im wesentlichen dazu beitragen, einen zuverlässigen und positiven Eindruck unserer Identität im Blockchain-Bereich durchzusetzen.", 2)

# 4) "Official font is " -> "Offizielle Schriftart ist "
$d.Content.Find.Execute("Official font is ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Offizielle Schriftart ist ", 2)
